$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 1082:1083, shifting the existing data
# (previously rows 1082-1178) down to rows 1084-1180.
$ws.Rows("1082:1083").Insert()

# New row 1082 - "Primera" quality record for date 45132
$ws.Cells.Item(1082,1).Value = 3
$ws.Cells.Item(1082,2).Value = "Femacal de La Calera"
$ws.Cells.Item(1082,3).Value = "Coquimbo"
$ws.Cells.Item(1082,4).Value = 45132
$ws.Cells.Item(1082,5).Value = 5
$ws.Cells.Item(1082,6).Value = 100112023
$ws.Cells.Item(1082,7).Value = "Brócoli"
$ws.Cells.Item(1082,8).Value = "Sin especificar"
$ws.Cells.Item(1082,9).Value = "Primera"
$ws.Cells.Item(1082,10).Value = 1900
$ws.Cells.Item(1082,11).Value = 700
$ws.Cells.Item(1082,12).Value = 700
$ws.Cells.Item(1082,13).Value = 700
$ws.Cells.Item(1082,14).Value = "$/unidad"
$ws.Cells.Item(1082,15).Value = "Provincia de Quillota"
$ws.Cells.Item(1082,16).Value = 700
$ws.Cells.Item(1082,17).Value = 1
$ws.Cells.Item(1082,18).Value = "Hortaliza"

# New row 1083 - "Segunda" quality record for date 45132
$ws.Cells.Item(1083,1).Value = 3
$ws.Cells.Item(1083,2).Value = "Femacal de La Calera"
$ws.Cells.Item(1083,3).Value = "Coquimbo"
$ws.Cells.Item(1083,4).Value = 45132
$ws.Cells.Item(1083,5).Value = 5
$ws.Cells.Item(1083,6).Value = 100112023
$ws.Cells.Item(1083,7).Value = "Brócoli"
$ws.Cells.Item(1083,8).Value = "Sin especificar"
$ws.Cells.Item(1083,9).Value = "Segunda"
$ws.Cells.Item(1083,10).Value = 1200
$ws.Cells.Item(1083,11).Value = 550
$ws.Cells.Item(1083,12).Value = 550
$ws.Cells.Item(1083,13).Value = 550
$ws.Cells.Item(1083,14).Value = "$/unidad"
$ws.Cells.Item(1083,15).Value = "Provincia de Quillota"
$ws.Cells.Item(1083,16).Value = 550
$ws.Cells.Item(1083,17).Value = 1
$ws.Cells.Item(1083,18).Value = "Hortaliza"
